$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bulk row-region updates (column A counters and column AN totals)
$ws.Range("A18:A281").Value = 6
$ws.Range("AN18:AN285").Value = 10

# Row 17 values bumped from 1 to 21, except AG17 which drops to 0
$ws.Range("B17:AF17").Value = 21
$ws.Range("AG17").Value = 0
$ws.Range("AH17:AM17").Value = 21

# Rows 18-24: AG column turns on (0 -> 1)
$ws.Range("AG18:AG24").Value = 1

# Row 29: redistributed sound-effect weighting across Z:AG
$ws.Range("Z29").Value = 5
$ws.Range("AA29:AB29").Value = 2
$ws.Range("AE29:AF29").Value = 2
$ws.Range("AG29").Value = 3

# Row 282: widen the active span and bump weight
$ws.Range("B282:O282").Value = 2
$ws.Range("P282").Value = 22

# Rows 283-285: activate column B, bump column P
$ws.Range("B283:B285").Value = 1
$ws.Range("P283:P285").Value = 6

# Row 286: bump Q:AM from 1 to 2
$ws.Range("Q286:AM286").Value = 2

# Window/view adjustments to match the saved workbook state
$excel.ActiveWindow.WindowState = -4143
$wb.Windows.Item(1).Left = 14295
$wb.Windows.Item(1).Top = 0
$wb.Windows.Item(1).Width = 14610
$wb.Windows.Item(1).Height = 16305

$ws.Activate()
$ws.Range("A49").Select()
$excel.ActiveWindow.ScrollRow = 49
$ws.Range("A1:AN300").Select()
$excel.Application.Goto($ws.Range("AN300"), $false)
